# feat: add command to attach tags to notes and vocabulary entries
#
# This particular export batch adds 15 freshly-studied vocabulary words
# (rows 81-95, captured 2021-11-14/2021-11-15) to the ENGLISH sheet. The
# NOTES sheet is untouched by this batch.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENGLISH")

# Each row: Word, Definition, Synonyms, Antonyms, Created at
# (Correct answer count is always 0 for a freshly-added entry.)
$entries = @(
    @("fatigue", "", "tiredness", "", "2021-11-14 20:15:38.955425"),
    @("cumulatively", "in a way that increases in quantity, degree, or force by successive addition", "", "", "2021-11-14 20:16:29.606448"),
    @("ample", "", "enough;plentiful", "", "2021-11-14 20:17:00.34091"),
    @("pernicious", "", "harmful", "", "2021-11-14 20:17:55.875123"),
    @("liken", "point out the resemblance of someone or smth to", "compare", "", "2021-11-14 20:18:50.405149"),
    @("machismo", "aggresive masculinity", "", "", "2021-11-15 14:49:54.205507"),
    @("stigma", "a mark of disgrace associated with a particular circumstance, quality, or person", "shame", "", "2021-11-15 14:50:45.332093"),
    @("commodity", "", "item", "", "2021-11-15 14:51:17.949615"),
    @("deride", "express contempt for", "ridicule", "", "2021-11-15 14:51:56.068144"),
    @("wimpish", "weak and cowardly or unadventurous", "", "", "2021-11-15 14:52:35.494874"),
    @("tout", "attempt to sell smth, typically by a direct or persistent approach", "peddle;solicit", "", "2021-11-15 14:53:45.152971"),
    @("discerning", "having or showing good judgement", "", "", "2021-11-15 14:54:30.922667"),
    @("mediocre", "of only average quality; not very good", "ordinary", "", "2021-11-15 14:55:50.807842"),
    @("utter", "", "absolute;complete", "", "2021-11-15 14:56:38.984503"),
    @("conviction", "", "belief", "", "2021-11-15 14:57:11.455154")
)

$startRow = 81
for ($i = 0; $i -lt $entries.Count; $i++) {
    $row = $startRow + $i
    $entry = $entries[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = 0
    $ws.Cells.Item($row, 6).Value = $entry[4]
}

Write-Host "Added $($entries.Count) vocabulary rows to ENGLISH sheet (now $($ws.UsedRange.Rows.Count) rows)"
